$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task list row that used to read "admin dashboard(Edit product)" was
# reworded to call out that it also covers the orders page.
$ws.Range("C18").Value = "admin dashboard(Edit product - orders page)"

# The author's last on-screen selection when the file was saved was C18.
$ws.Range("C18").Select()
